# act 08-08-2025 se agregan mas equipos y nueva funcionalidad
# Fill in match-result data for row 2 (Karlsruher vs Munster, 02/08/2025)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fecha (A2) - force text so Excel doesn't auto-convert the dd/mm/yyyy string
# into a date serial number, then restore the default "Normal" style so no
# extra formatting sticks to the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "02/08/2025"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "Karlsruher"      # Equipo Local
$ws.Range("C2").Value = 3                 # Goles Local
$ws.Range("D2").Value = 2                 # Goles Visitante
$ws.Range("E2").Value = "Munster"         # Visitante
$ws.Range("F2").Value = "L"               # Resultado
$ws.Range("G2").Value = 1                 # 1T Goles Favor
$ws.Range("H2").Value = 2                 # 1T Goles Contra
$ws.Range("I2").Value = 1                 # 2T Goles Favor
$ws.Range("J2").Value = 1                 # 2T Goles Contra
$ws.Range("K2").Value = 0.64              # xG Favor
$ws.Range("L2").Value = 1.68              # xG Contra
$ws.Range("M2").Value = 11                # Shots Favor
$ws.Range("N2").Value = 13                # Shots Contra
$ws.Range("O2").Value = 3                 # A puerta Favor
$ws.Range("P2").Value = 4                 # A puerta Contra
